# Apply updated crypto price/volume data to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.104.61"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "2.220.62"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.96%  "
$ws.Range("E6").Value = "  -1.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.15"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0962"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.17"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").Value = "2.552.84"
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.847"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("D17").Value = "2.217.79"
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("D18").Value = "41.904.83"
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000109"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +11.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +33.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.96%  "
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("E28").Value = "  -1.20%  "
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +14.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0802"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.47%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.114"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.48%  "
$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "29.30"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.90%  "
$ws.Range("E37").Value = "  -5.22%  "
$ws.Range("E38").Value = "  -1.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.05"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.07%  "
$ws.Range("E40").Value = "  -1.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "65.47"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.17%  "
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("E44").Value = "  +1.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.89%  "
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.55%  "
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("D51").Value = "2.427.94"
$ws.Range("E51").Value = "  -1.27%  "
